# ---------------------------------------------------------------
# hl60 RNA-seq 72h used instead of 120h
#
# 1) "Comparisons" sheet: drop the six Stbl.* (stability-assay)
#    column-pairs (N:Y) -- they shift the trailing TE.Estimate /
#    TE.fdr pair left into N:O -- and refresh the hl60 Exp values
#    (+ matching row order) to the 72h re-analysis.
# 2) "CRISPRi-screens" sheet: re-sort rows to the same gene order.
# ---------------------------------------------------------------
$wb = $excel.ActiveWorkbook

# --- Comparisons sheet -----------------------------------------
$ws2 = $wb.Worksheets.Item("Comparisons")

# Remove the Stbl.hl60 / Stbl.kg1 / Stbl.molm14 / Stbl.ociaml2 /
# Stbl.ociaml3 / Stbl.thp1 column-pairs; this shifts the old Z:AA
# (TE.Estimate_treatmentDRUG / TE.fdr_Pr...z.._treatmentDRUG) pair
# left so it lands in N:O.
$ws2.Range("N1:Y7").EntireColumn.Delete()

# Rewrite rows 2:7 in the updated gene order with the refreshed
# hl60 72h values (logFC for several genes changed; pvals and all
# other assay columns are unchanged).
$ws2.Cells.Item(2,1).Value = "SHOC2"
$ws2.Cells.Item(2,2).Value = -0.8149999999999999
$ws2.Cells.Item(2,3).Value = 0.001
$ws2.Cells.Item(2,4).Value = -0.038
$ws2.Cells.Item(2,5).Value = 0.852
$ws2.Cells.Item(2,6).Value = -0.194
$ws2.Cells.Item(2,7).Value = 0.276
$ws2.Cells.Item(2,8).Value = -0.635
$ws2.Cells.Item(2,9).Value = 0.037
$ws2.Cells.Item(2,10).Value = 0.02
$ws2.Cells.Item(2,11).Value = 0.914
$ws2.Cells.Item(2,12).Value = -0.182
$ws2.Cells.Item(2,13).Value = 0.359
$ws2.Cells.Item(2,14).Value = -0.401
$ws2.Cells.Item(2,15).Value = 0.5639999999999999

$ws2.Cells.Item(3,1).Value = "DENND1A"
$ws2.Cells.Item(3,2).Value = -0.344
$ws2.Cells.Item(3,3).Value = 0.018
$ws2.Cells.Item(3,4).Value = -0.062
$ws2.Cells.Item(3,5).Value = 0.8100000000000001
$ws2.Cells.Item(3,6).Value = -0.334
$ws2.Cells.Item(3,7).Value = 0.166
$ws2.Cells.Item(3,8).Value = -0.037
$ws2.Cells.Item(3,9).Value = 0.89
$ws2.Cells.Item(3,10).Value = -0.408
$ws2.Cells.Item(3,11).Value = 0.102
$ws2.Cells.Item(3,12).Value = -0.762
$ws2.Cells.Item(3,13).Value = 0
$ws2.Cells.Item(3,14).Value = 0.474
$ws2.Cells.Item(3,15).Value = 0.6889999999999999

$ws2.Cells.Item(4,1).Value = "ZNF830"
$ws2.Cells.Item(4,2).Value = -0.314
$ws2.Cells.Item(4,3).Value = 0.006
$ws2.Cells.Item(4,4).Value = -0.455
$ws2.Cells.Item(4,5).Value = 0.033
$ws2.Cells.Item(4,6).Value = -0.073
$ws2.Cells.Item(4,7).Value = 0.779
$ws2.Cells.Item(4,8).Value = 0.222
$ws2.Cells.Item(4,9).Value = 0.415
$ws2.Cells.Item(4,10).Value = -0.203
$ws2.Cells.Item(4,11).Value = 0.444
$ws2.Cells.Item(4,12).Value = -0.336
$ws2.Cells.Item(4,13).Value = 0.129
$ws2.Cells.Item(4,14).Value = 0.116
$ws2.Cells.Item(4,15).Value = 0.945

$ws2.Cells.Item(5,1).Value = "ZC3H13"
$ws2.Cells.Item(5,2).Value = -1.391
$ws2.Cells.Item(5,3).Value = 0.01
$ws2.Cells.Item(5,4).Value = 0.199
$ws2.Cells.Item(5,5).Value = 0.514
$ws2.Cells.Item(5,6).Value = -0.064
$ws2.Cells.Item(5,7).Value = 0.765
$ws2.Cells.Item(5,8).Value = -0.006
$ws2.Cells.Item(5,9).Value = 0.981
$ws2.Cells.Item(5,10).Value = 0.38
$ws2.Cells.Item(5,11).Value = 0.194
$ws2.Cells.Item(5,12).Value = -0.154
$ws2.Cells.Item(5,13).Value = 0.516
$ws2.Cells.Item(5,14).Value = 0.077
$ws2.Cells.Item(5,15).Value = 0.93

$ws2.Cells.Item(6,1).Value = "MANBA"
$ws2.Cells.Item(6,2).Value = -0.347
$ws2.Cells.Item(6,3).Value = 0.006
$ws2.Cells.Item(6,4).Value = -0.004
$ws2.Cells.Item(6,5).Value = 0.982
$ws2.Cells.Item(6,6).Value = 0.298
$ws2.Cells.Item(6,7).Value = 0.245
$ws2.Cells.Item(6,8).Value = 0.064
$ws2.Cells.Item(6,9).Value = 0.784
$ws2.Cells.Item(6,10).Value = 0.036
$ws2.Cells.Item(6,11).Value = 0.903
$ws2.Cells.Item(6,12).Value = -0.182
$ws2.Cells.Item(6,13).Value = 0.277
$ws2.Cells.Item(6,14).Value = 0.219
$ws2.Cells.Item(6,15).Value = 0.751

$ws2.Cells.Item(7,1).Value = "PTCD1"
$ws2.Cells.Item(7,2).Value = -0.213
$ws2.Cells.Item(7,3).Value = 0
$ws2.Cells.Item(7,4).Value = 1.101
$ws2.Cells.Item(7,5).Value = 0.002
$ws2.Cells.Item(7,6).Value = -0.113
$ws2.Cells.Item(7,7).Value = 0.785
$ws2.Cells.Item(7,8).Value = -0.249
$ws2.Cells.Item(7,9).Value = 0.488
$ws2.Cells.Item(7,10).Value = -0.251
$ws2.Cells.Item(7,11).Value = 0.5
$ws2.Cells.Item(7,12).Value = 0.073
$ws2.Cells.Item(7,13).Value = 0.884
$ws2.Cells.Item(7,14).ClearContents()
$ws2.Cells.Item(7,15).ClearContents()

# --- CRISPRi-screens sheet ---------------------------------------
# Same six genes, same values -- only the row order changes (to
# line up with the "Comparisons" sheet above).
$ws3 = $wb.Worksheets.Item("CRISPRi-screens")
$ws3.Cells.Item(2,1).Value = "SHOC2"
$ws3.Cells.Item(2,2).Value = 0.603
$ws3.Cells.Item(2,3).Value = 0
$ws3.Cells.Item(2,4).Value = 0.365
$ws3.Cells.Item(2,5).Value = 0.08699999999999999
$ws3.Cells.Item(2,6).Value = 0.676
$ws3.Cells.Item(2,7).Value = 0.494
$ws3.Cells.Item(2,8).Value = -0.351
$ws3.Cells.Item(2,9).Value = 0
$ws3.Cells.Item(2,10).Value = -0.644
$ws3.Cells.Item(2,11).Value = 0.001
$ws3.Cells.Item(2,12).Value = -0.644
$ws3.Cells.Item(2,13).Value = 0.001

$ws3.Cells.Item(3,1).Value = "DENND1A"
$ws3.Cells.Item(3,2).Value = 0.104
$ws3.Cells.Item(3,3).Value = 0.546
$ws3.Cells.Item(3,4).Value = 0.112
$ws3.Cells.Item(3,5).Value = 0.158
$ws3.Cells.Item(3,6).Value = 0.092
$ws3.Cells.Item(3,7).Value = 0.304
$ws3.Cells.Item(3,8).Value = -0.057
$ws3.Cells.Item(3,9).Value = 0.977
$ws3.Cells.Item(3,10).Value = -0.022
$ws3.Cells.Item(3,11).Value = 0.892
$ws3.Cells.Item(3,12).Value = -0.022
$ws3.Cells.Item(3,13).Value = 0.892

$ws3.Cells.Item(4,1).Value = "ZNF830"
$ws3.Cells.Item(4,2).Value = 0.129
$ws3.Cells.Item(4,3).Value = 0.096
$ws3.Cells.Item(4,4).Value = 0.203
$ws3.Cells.Item(4,5).Value = 0.226
$ws3.Cells.Item(4,6).Value = 0.048
$ws3.Cells.Item(4,7).Value = 0.51
$ws3.Cells.Item(4,8).Value = -0.277
$ws3.Cells.Item(4,9).Value = 0.042
$ws3.Cells.Item(4,10).Value = -0.09
$ws3.Cells.Item(4,11).Value = 0.128
$ws3.Cells.Item(4,12).Value = -0.09
$ws3.Cells.Item(4,13).Value = 0.128

$ws3.Cells.Item(5,1).Value = "ZC3H13"
$ws3.Cells.Item(5,2).Value = 0.283
$ws3.Cells.Item(5,3).Value = 0
$ws3.Cells.Item(5,4).Value = 0.216
$ws3.Cells.Item(5,5).Value = 0.058
$ws3.Cells.Item(5,6).Value = 0.379
$ws3.Cells.Item(5,7).Value = 0.005
$ws3.Cells.Item(5,8).Value = -0.102
$ws3.Cells.Item(5,9).Value = 0.001
$ws3.Cells.Item(5,10).Value = -0.156
$ws3.Cells.Item(5,11).Value = 0.046
$ws3.Cells.Item(5,12).Value = -0.156
$ws3.Cells.Item(5,13).Value = 0.046

$ws3.Cells.Item(6,1).Value = "MANBA"
$ws3.Cells.Item(6,2).Value = 0.103
$ws3.Cells.Item(6,3).Value = 0.022
$ws3.Cells.Item(6,4).Value = 0.184
$ws3.Cells.Item(6,5).Value = 0.016
$ws3.Cells.Item(6,6).Value = 0.163
$ws3.Cells.Item(6,7).Value = 0.138
$ws3.Cells.Item(6,8).Value = -0.014
$ws3.Cells.Item(6,9).Value = 0.371
$ws3.Cells.Item(6,10).Value = -0.056
$ws3.Cells.Item(6,11).Value = 0.6889999999999999
$ws3.Cells.Item(6,12).Value = -0.056
$ws3.Cells.Item(6,13).Value = 0.6889999999999999

$ws3.Cells.Item(7,1).Value = "PTCD1"
$ws3.Cells.Item(7,2).Value = 0.16
$ws3.Cells.Item(7,3).Value = 0.099
$ws3.Cells.Item(7,4).Value = 0.194
$ws3.Cells.Item(7,5).Value = 0.023
$ws3.Cells.Item(7,6).Value = 0.166
$ws3.Cells.Item(7,7).Value = 0.028
$ws3.Cells.Item(7,8).Value = -0.083
$ws3.Cells.Item(7,9).Value = 0.015
$ws3.Cells.Item(7,10).Value = -0.198
$ws3.Cells.Item(7,11).Value = 0.001
$ws3.Cells.Item(7,12).Value = -0.198
$ws3.Cells.Item(7,13).Value = 0.001

